# fix: pt2 - excluding contratada material nosense
# Adjusts quantity (col F) and total value (col H) for a handful of
# stock rows on Sheet1 to remove "contratada" material that shouldn't
# have been counted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43
$ws.Range("F43").Value = 107
$ws.Range("H43").Value = 16964.31

# Row 44
$ws.Range("F44").Value = 26
$ws.Range("H44").Value = 14528.15

# Row 45
$ws.Range("F45").Value = 58
$ws.Range("H45").Value = 5113.09

# Row 55
$ws.Range("F55").Value = 560
$ws.Range("H55").Value = 5081.82

# Row 57
$ws.Range("F57").Value = 3050
$ws.Range("H57").Value = 31347.92

# Row 81 (quantity unchanged, only rounding of the total value)
$ws.Range("H81").Value = 8076.41

# Row 83
$ws.Range("F83").Value = 2796
$ws.Range("H83").Value = 72828.85

# Row 98
$ws.Range("F98").Value = 2365
$ws.Range("H98").Value = 332579.30

# Row 121
$ws.Range("F121").Value = 12271
$ws.Range("H121").Value = 47632.80

# Row 174
$ws.Range("F174").Value = 2441
$ws.Range("H174").Value = 45707.72

# Row 188 (quantity unchanged, only rounding of the total value)
$ws.Range("H188").Value = 7253.54

# Row 230
$ws.Range("F230").Value = 20292.8
$ws.Range("H230").Value = 46050.26

# Row 245
$ws.Range("F245").Value = 5601
$ws.Range("H245").Value = 14060.56
